$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "C3"
$ws.Cells.Item(2, 3).Value = "Itgam"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 3.820425
$ws.Cells.Item(2, 8).Value = 11.461275
$ws.Cells.Item(2, 9).Value = 0.02049663039797357
$ws.Cells.Item(2, 10).Value = 0.02049663039797357
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.142723
$ws.Cells.Item(2, 14).Value = 0.428169
$ws.Cells.Item(2, 15).Value = 0.0009642800942465787
$ws.Cells.Item(2, 16).Value = 0.0009642800942465787
$ws.Cells.Item(2, 17).Value = 0.5452625172750001
$ws.Cells.Item(2, 18).Value = 4.907362655475001
$ws.Cells.Item(2, 19).Value = [double]"1.976449269189524E-05"
$ws.Cells.Item(2, 20).Value = [double]"1.976449269189524E-05"

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "C3"
$ws.Cells.Item(3, 3).Value = "Itgam"
$ws.Cells.Item(3, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 3.820425
$ws.Cells.Item(3, 8).Value = 11.461275
$ws.Cells.Item(3, 9).Value = 0.02049663039797357
$ws.Cells.Item(3, 10).Value = 0.02049663039797357
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 86.42780700000002
$ws.Cells.Item(3, 14).Value = 259.283421
$ws.Cells.Item(3, 15).Value = 0.5839326098770704
$ws.Cells.Item(3, 16).Value = 0.5839326098770704
$ws.Cells.Item(3, 17).Value = 330.1909545579751
$ws.Cells.Item(3, 18).Value = 2971.718591021775
$ws.Cells.Item(3, 19).Value = 0.0119686508819744
$ws.Cells.Item(3, 20).Value = 0.0119686508819744

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "C3"
$ws.Cells.Item(4, 3).Value = "Itgam"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 3.820425
$ws.Cells.Item(4, 8).Value = 11.461275
$ws.Cells.Item(4, 9).Value = 0.02049663039797357
$ws.Cells.Item(4, 10).Value = 0.02049663039797357
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.006361333333333333
$ws.Cells.Item(4, 14).Value = 0.019084
$ws.Cells.Item(4, 15).Value = [double]"4.297910712499435E-05"
$ws.Cells.Item(4, 16).Value = [double]"4.297910712499435E-05"
$ws.Cells.Item(4, 17).Value = 0.0243029969
$ws.Cells.Item(4, 18).Value = 0.2187269721
$ws.Cells.Item(4, 19).Value = [double]"8.809268735759215E-07"
$ws.Cells.Item(4, 20).Value = [double]"8.809268735759215E-07"

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "C3"
$ws.Cells.Item(5, 3).Value = "Itgam"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 3.820425
$ws.Cells.Item(5, 8).Value = 11.461275
$ws.Cells.Item(5, 9).Value = 0.02049663039797357
$ws.Cells.Item(5, 10).Value = 0.02049663039797357
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 61.43300833333333
$ws.Cells.Item(5, 14).Value = 184.299025
$ws.Cells.Item(5, 15).Value = 0.415060130921558
$ws.Cells.Item(5, 16).Value = 0.415060130921558
$ws.Cells.Item(5, 17).Value = 234.700200861875
$ws.Cells.Item(5, 18).Value = 2112.301807756875
$ws.Cells.Item(5, 19).Value = 0.008507334096433693
$ws.Cells.Item(5, 20).Value = 0.008507334096433693

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "C3"
$ws.Cells.Item(6, 3).Value = "Itgam"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 140.5890806666667
$ws.Cells.Item(6, 8).Value = 421.767242
$ws.Cells.Item(6, 9).Value = 0.7542622677884155
$ws.Cells.Item(6, 10).Value = 0.7542622677884157
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.142723
$ws.Cells.Item(6, 14).Value = 0.428169
$ws.Cells.Item(6, 15).Value = 0.0009642800942465787
$ws.Cells.Item(6, 16).Value = 0.0009642800942465787
$ws.Cells.Item(6, 17).Value = 20.06529535998867
$ws.Cells.Item(6, 18).Value = 180.587658239898
$ws.Cells.Item(6, 19).Value = 0.0007273200906696515
$ws.Cells.Item(6, 20).Value = 0.0007273200906696516

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "C3"
$ws.Cells.Item(7, 3).Value = "Itgam"
$ws.Cells.Item(7, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 140.5890806666667
$ws.Cells.Item(7, 8).Value = 421.767242
$ws.Cells.Item(7, 9).Value = 0.7542622677884155
$ws.Cells.Item(7, 10).Value = 0.7542622677884157
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 86.42780700000002
$ws.Cells.Item(7, 14).Value = 259.283421
$ws.Cells.Item(7, 15).Value = 0.5839326098770704
$ws.Cells.Item(7, 16).Value = 0.5839326098770704
$ws.Cells.Item(7, 17).Value = 12150.8059301661
$ws.Cells.Item(7, 18).Value = 109357.2533714949
$ws.Cells.Item(7, 19).Value = 0.4404383345614873
$ws.Cells.Item(7, 20).Value = 0.4404383345614873

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "C3"
$ws.Cells.Item(8, 3).Value = "Itgam"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 140.5890806666667
$ws.Cells.Item(8, 8).Value = 421.767242
$ws.Cells.Item(8, 9).Value = 0.7542622677884155
$ws.Cells.Item(8, 10).Value = 0.7542622677884157
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.006361333333333333
$ws.Cells.Item(8, 14).Value = 0.019084
$ws.Cells.Item(8, 15).Value = [double]"4.297910712499435E-05"
$ws.Cells.Item(8, 16).Value = [double]"4.297910712499435E-05"
$ws.Cells.Item(8, 17).Value = 0.8943340051475556
$ws.Cells.Item(8, 18).Value = 8.049006046328
$ws.Cells.Item(8, 19).Value = [double]"3.241751880761949E-05"
$ws.Cells.Item(8, 20).Value = [double]"3.241751880761949E-05"

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "C3"
$ws.Cells.Item(9, 3).Value = "Itgam"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 140.5890806666667
$ws.Cells.Item(9, 8).Value = 421.767242
$ws.Cells.Item(9, 9).Value = 0.7542622677884155
$ws.Cells.Item(9, 10).Value = 0.7542622677884157
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 61.43300833333333
$ws.Cells.Item(9, 14).Value = 184.299025
$ws.Cells.Item(9, 15).Value = 0.415060130921558
$ws.Cells.Item(9, 16).Value = 0.415060130921558
$ws.Cells.Item(9, 17).Value = 8636.810164171005
$ws.Cells.Item(9, 18).Value = 77731.29147753905
$ws.Cells.Item(9, 19).Value = 0.313064195617451
$ws.Cells.Item(9, 20).Value = 0.313064195617451

# Row 10
$ws.Cells.Item(10, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 2).Value = "C3"
$ws.Cells.Item(10, 3).Value = "Itgam"
$ws.Cells.Item(10, 4).Value = "FAPs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 30.51067
$ws.Cells.Item(10, 8).Value = 91.53201
$ws.Cells.Item(10, 9).Value = 0.1636901460399144
$ws.Cells.Item(10, 10).Value = 0.1636901460399144
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.142723
$ws.Cells.Item(10, 14).Value = 0.428169
$ws.Cells.Item(10, 15).Value = 0.0009642800942465787
$ws.Cells.Item(10, 16).Value = 0.0009642800942465787
$ws.Cells.Item(10, 17).Value = 4.35457435441
$ws.Cells.Item(10, 18).Value = 39.19116918969
$ws.Cells.Item(10, 19).Value = 0.0001578431494506049
$ws.Cells.Item(10, 20).Value = 0.0001578431494506049

# Row 11
$ws.Cells.Item(11, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(11, 2).Value = "C3"
$ws.Cells.Item(11, 3).Value = "Itgam"
$ws.Cells.Item(11, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 30.51067
$ws.Cells.Item(11, 8).Value = 91.53201
$ws.Cells.Item(11, 9).Value = 0.1636901460399144
$ws.Cells.Item(11, 10).Value = 0.1636901460399144
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 86.42780700000002
$ws.Cells.Item(11, 14).Value = 259.283421
$ws.Cells.Item(11, 15).Value = 0.5839326098770704
$ws.Cells.Item(11, 16).Value = 0.5839326098770704
$ws.Cells.Item(11, 17).Value = 2636.970298200691
$ws.Cells.Item(11, 18).Value = 23732.73268380621
$ws.Cells.Item(11, 19).Value = 0.09558401418824604
$ws.Cells.Item(11, 20).Value = 0.09558401418824604

# Row 12
$ws.Cells.Item(12, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(12, 2).Value = "C3"
$ws.Cells.Item(12, 3).Value = "Itgam"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 30.51067
$ws.Cells.Item(12, 8).Value = 91.53201
$ws.Cells.Item(12, 9).Value = 0.1636901460399144
$ws.Cells.Item(12, 10).Value = 0.1636901460399144
$ws.Cells.Item(12, 11).Value = 1
$ws.Cells.Item(12, 12).Value = 0.3333333333333333
$ws.Cells.Item(12, 13).Value = 0.006361333333333333
$ws.Cells.Item(12, 14).Value = 0.019084
$ws.Cells.Item(12, 15).Value = [double]"4.297910712499435E-05"
$ws.Cells.Item(12, 16).Value = [double]"4.297910712499435E-05"
$ws.Cells.Item(12, 17).Value = 0.1940885420933333
$ws.Cells.Item(12, 18).Value = 1.74679687884
$ws.Cells.Item(12, 19).Value = [double]"7.035256321955453E-06"
$ws.Cells.Item(12, 20).Value = [double]"7.035256321955453E-06"

# Row 13
$ws.Cells.Item(13, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(13, 2).Value = "C3"
$ws.Cells.Item(13, 3).Value = "Itgam"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 30.51067
$ws.Cells.Item(13, 8).Value = 91.53201
$ws.Cells.Item(13, 9).Value = 0.1636901460399144
$ws.Cells.Item(13, 10).Value = 0.1636901460399144
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 61.43300833333333
$ws.Cells.Item(13, 14).Value = 184.299025
$ws.Cells.Item(13, 15).Value = 0.415060130921558
$ws.Cells.Item(13, 16).Value = 0.415060130921558
$ws.Cells.Item(13, 17).Value = 1874.362244365583
$ws.Cells.Item(13, 18).Value = 16869.26019929025
$ws.Cells.Item(13, 19).Value = 0.06794125344589583
$ws.Cells.Item(13, 20).Value = 0.06794125344589583

# Row 14
$ws.Cells.Item(14, 1).Value = "MuSCs"
$ws.Cells.Item(14, 2).Value = "C3"
$ws.Cells.Item(14, 3).Value = "Itgam"
$ws.Cells.Item(14, 4).Value = "FAPs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.258813
$ws.Cells.Item(14, 8).Value = 0.776439
$ws.Cells.Item(14, 9).Value = 0.001388535150720334
$ws.Cells.Item(14, 10).Value = 0.001388535150720334
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.142723
$ws.Cells.Item(14, 14).Value = 0.428169
$ws.Cells.Item(14, 15).Value = 0.0009642800942465787
$ws.Cells.Item(14, 16).Value = 0.0009642800942465787
$ws.Cells.Item(14, 17).Value = 0.036938567799
$ws.Cells.Item(14, 18).Value = 0.332447110191
$ws.Cells.Item(14, 19).Value = [double]"1.338936806001291E-06"
$ws.Cells.Item(14, 20).Value = [double]"1.338936806001291E-06"

# Row 15
$ws.Cells.Item(15, 1).Value = "MuSCs"
$ws.Cells.Item(15, 2).Value = "C3"
$ws.Cells.Item(15, 3).Value = "Itgam"
$ws.Cells.Item(15, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.258813
$ws.Cells.Item(15, 8).Value = 0.776439
$ws.Cells.Item(15, 9).Value = 0.001388535150720334
$ws.Cells.Item(15, 10).Value = 0.001388535150720334
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 86.42780700000002
$ws.Cells.Item(15, 14).Value = 259.283421
$ws.Cells.Item(15, 15).Value = 0.5839326098770704
$ws.Cells.Item(15, 16).Value = 0.5839326098770704
$ws.Cells.Item(15, 17).Value = 22.36864001309101
$ws.Cells.Item(15, 18).Value = 201.317760117819
$ws.Cells.Item(15, 19).Value = 0.000810810954466176
$ws.Cells.Item(15, 20).Value = 0.0008108109544661759

# Row 16
$ws.Cells.Item(16, 1).Value = "MuSCs"
$ws.Cells.Item(16, 2).Value = "C3"
$ws.Cells.Item(16, 3).Value = "Itgam"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.258813
$ws.Cells.Item(16, 8).Value = 0.776439
$ws.Cells.Item(16, 9).Value = 0.001388535150720334
$ws.Cells.Item(16, 10).Value = 0.001388535150720334
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.006361333333333333
$ws.Cells.Item(16, 14).Value = 0.019084
$ws.Cells.Item(16, 15).Value = [double]"4.297910712499435E-05"
$ws.Cells.Item(16, 16).Value = [double]"4.297910712499435E-05"
$ws.Cells.Item(16, 17).Value = 0.001646395764
$ws.Cells.Item(16, 18).Value = 0.014817561876
$ws.Cells.Item(16, 19).Value = [double]"5.967800098962942E-08"
$ws.Cells.Item(16, 20).Value = [double]"5.967800098962941E-08"

# Row 17
$ws.Cells.Item(17, 1).Value = "MuSCs"
$ws.Cells.Item(17, 2).Value = "C3"
$ws.Cells.Item(17, 3).Value = "Itgam"
$ws.Cells.Item(17, 4).Value = "Resolving-Mac"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.258813
$ws.Cells.Item(17, 8).Value = 0.776439
$ws.Cells.Item(17, 9).Value = 0.001388535150720334
$ws.Cells.Item(17, 10).Value = 0.001388535150720334
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 61.43300833333333
$ws.Cells.Item(17, 14).Value = 184.299025
$ws.Cells.Item(17, 15).Value = 0.415060130921558
$ws.Cells.Item(17, 16).Value = 0.415060130921558
$ws.Cells.Item(17, 17).Value = 15.899661185775
$ws.Cells.Item(17, 18).Value = 143.096950671975
$ws.Cells.Item(17, 19).Value = 0.0005763255814471671
$ws.Cells.Item(17, 20).Value = 0.000576325581447167

# Row 18
$ws.Cells.Item(18, 1).Value = "Resolving-Mac"
$ws.Cells.Item(18, 2).Value = "C3"
$ws.Cells.Item(18, 3).Value = "Itgam"
$ws.Cells.Item(18, 4).Value = "FAPs"
$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 11.213844
$ws.Cells.Item(18, 8).Value = 33.641532
$ws.Cells.Item(18, 9).Value = 0.0601624206229761
$ws.Cells.Item(18, 10).Value = 0.0601624206229761
$ws.Cells.Item(18, 11).Value = 1
$ws.Cells.Item(18, 12).Value = 0.3333333333333333
$ws.Cells.Item(18, 13).Value = 0.142723
$ws.Cells.Item(18, 14).Value = 0.428169
$ws.Cells.Item(18, 15).Value = 0.0009642800942465787
$ws.Cells.Item(18, 16).Value = 0.0009642800942465787
$ws.Cells.Item(18, 17).Value = 1.600473457212
$ws.Cells.Item(18, 18).Value = 14.404261114908
$ws.Cells.Item(18, 19).Value = [double]"5.801342462842571E-05"
$ws.Cells.Item(18, 20).Value = [double]"5.80134246284257E-05"

# Row 19
$ws.Cells.Item(19, 1).Value = "Resolving-Mac"
$ws.Cells.Item(19, 2).Value = "C3"
$ws.Cells.Item(19, 3).Value = "Itgam"
$ws.Cells.Item(19, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 11.213844
$ws.Cells.Item(19, 8).Value = 33.641532
$ws.Cells.Item(19, 9).Value = 0.0601624206229761
$ws.Cells.Item(19, 10).Value = 0.0601624206229761
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 86.42780700000002
$ws.Cells.Item(19, 14).Value = 259.283421
$ws.Cells.Item(19, 15).Value = 0.5839326098770704
$ws.Cells.Item(19, 16).Value = 0.5839326098770704
$ws.Cells.Item(19, 17).Value = 969.1879449601082
$ws.Cells.Item(19, 18).Value = 8722.691504640972
$ws.Cells.Item(19, 19).Value = 0.03513079929089652
$ws.Cells.Item(19, 20).Value = 0.03513079929089651

# Row 20
$ws.Cells.Item(20, 1).Value = "Resolving-Mac"
$ws.Cells.Item(20, 2).Value = "C3"
$ws.Cells.Item(20, 3).Value = "Itgam"
$ws.Cells.Item(20, 4).Value = "MuSCs"
$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 11.213844
$ws.Cells.Item(20, 8).Value = 33.641532
$ws.Cells.Item(20, 9).Value = 0.0601624206229761
$ws.Cells.Item(20, 10).Value = 0.0601624206229761
$ws.Cells.Item(20, 11).Value = 1
$ws.Cells.Item(20, 12).Value = 0.3333333333333333
$ws.Cells.Item(20, 13).Value = 0.006361333333333333
$ws.Cells.Item(20, 14).Value = 0.019084
$ws.Cells.Item(20, 15).Value = [double]"4.297910712499435E-05"
$ws.Cells.Item(20, 16).Value = [double]"4.297910712499435E-05"
$ws.Cells.Item(20, 17).Value = 0.071334999632
$ws.Cells.Item(20, 18).Value = 0.642014996688
$ws.Cells.Item(20, 19).Value = [double]"2.585727120853859E-06"
$ws.Cells.Item(20, 20).Value = [double]"2.585727120853859E-06"

# Row 21
$ws.Cells.Item(21, 1).Value = "Resolving-Mac"
$ws.Cells.Item(21, 2).Value = "C3"
$ws.Cells.Item(21, 3).Value = "Itgam"
$ws.Cells.Item(21, 4).Value = "Resolving-Mac"
$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = 11.213844
$ws.Cells.Item(21, 8).Value = 33.641532
$ws.Cells.Item(21, 9).Value = 0.0601624206229761
$ws.Cells.Item(21, 10).Value = 0.0601624206229761
$ws.Cells.Item(21, 11).Value = 3
$ws.Cells.Item(21, 12).Value = 1
$ws.Cells.Item(21, 13).Value = 61.43300833333333
$ws.Cells.Item(21, 14).Value = 184.299025
$ws.Cells.Item(21, 15).Value = 0.415060130921558
$ws.Cells.Item(21, 16).Value = 0.415060130921558
$ws.Cells.Item(21, 17).Value = 688.9001719007
$ws.Cells.Item(21, 18).Value = 6200.101547106299
$ws.Cells.Item(21, 19).Value = 0.0249710221803303
$ws.Cells.Item(21, 20).Value = 0.0249710221803303

